$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-11 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(6, 2, 5, 7),
    @(3, 3, 5, 5),
    @(2, 4, 5, 7),
    @(7, 4, 12, 13),
    @(1, 5, 5, 5),
    @(6, 5, 10, 12),
    @(4, 6, 5, 6),
    @(5, 6, 11, 11),
    @(6, 6, 16, 18),
    @(6, 6, 23, 25)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
